$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (sheet tab name encodes "through" date)
$ws.Name = "Through 2021-10-17"

# Row 11 (September) - 2021 columns only
$ws.Range("T11").Value = 5
$ws.Range("U11").Value = 173
$ws.Range("V11").Value = 0.0281

# Row 12 (October) - label + several columns
$ws.Range("A12").Value = "October (through 10-17)"
$ws.Range("C12").Value = 13
$ws.Range("D12").Value = 0.0714
$ws.Range("L12").Value = 43
$ws.Range("M12").Value = 0.0444
$ws.Range("R12").Value = 82
$ws.Range("U12").Value = 105

# Row 13 (Total)
$ws.Range("C13").Value = 209
$ws.Range("D13").Value = 0.1292
$ws.Range("L13").Value = 530
$ws.Range("M13").Value = 0.1062
$ws.Range("R13").Value = 930
$ws.Range("S13").Value = 0.0539
$ws.Range("T13").Value = 83
$ws.Range("U13").Value = 1271
$ws.Range("V13").Value = 0.0613
